# Insert a new price record as row 43 on the (single) data sheet,
# pushing the existing rows 43:118 down to 44:119.
#
# "Fruta / hortaliza, semanal" — weekly refresh of the Bruselas (repollito)
# price series for "Vega Modelo de Temuco": a new observation is prepended
# to the daily log, shifting every subsequent record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 43:118 -> 44:119 (dimension grows from R118 to R119 automatically).
$ws.Rows("43:43").Insert()

# Populate the freshly inserted row 43 with the new weekly observation.
$ws.Cells.Item(43, 1).Value  = 10
$ws.Cells.Item(43, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value  = "La Araucanía"
$ws.Cells.Item(43, 4).Value  = 44797
$ws.Cells.Item(43, 5).Value  = 9
$ws.Cells.Item(43, 6).Value  = 100112035
$ws.Cells.Item(43, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(43, 8).Value  = "Sin especificar"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 80
$ws.Cells.Item(43, 11).Value = 24000
$ws.Cells.Item(43, 12).Value = 24000
$ws.Cells.Item(43, 13).Value = 24000
$ws.Cells.Item(43, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 2400
$ws.Cells.Item(43, 17).Value = 10
$ws.Cells.Item(43, 18).Value = "Hortaliza"
